# Apply the cryptos list refresh (GitHub Actions data update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range('D2')
$dCell.Value = '''63.049.92'
$dCell.Style = "Normal"
$ws.Range('E2').Value = '  +6.04%  '

$dCell = $ws.Range('D3')
$dCell.Value = '''3.112.61'
$dCell.Style = "Normal"
$ws.Range('E3').Value = '  +3.98%  '

$ws.Range('E4').Value = '  +0.00%  '

$dCell = $ws.Range('D5')
$dCell.Value = '''584.98'
$dCell.Style = "Normal"
$ws.Range('E5').Value = '  +3.91%  '

$dCell = $ws.Range('D6')
$dCell.Value = '''144.45'
$dCell.Style = "Normal"
$ws.Range('E6').Value = '  +3.77%  '

$ws.Range('E7').Value = '  +0.00%  '

$dCell = $ws.Range('D8')
$dCell.Value = '''3.104.39'
$dCell.Style = "Normal"
$ws.Range('E8').Value = '  +4.11%  '

$ws.Range('E9').Value = '  +1.64%  '

$dCell = $ws.Range('D10')
$dCell.Value = '''0.151'
$dCell.Style = "Normal"
$ws.Range('E10').Value = '  +13.34%  '

$dCell = $ws.Range('D11')
$dCell.Value = '''5.77'
$dCell.Style = "Normal"
$ws.Range('E11').Value = '  +8.97%  '

$dCell = $ws.Range('D12')
$dCell.Value = '''0.466'
$dCell.Style = "Normal"
$ws.Range('E12').Value = '  +2.91%  '

$ws.Range('E13').Value = '  +7.96%  '

$dCell = $ws.Range('D14')
$dCell.Value = '''35.53'
$dCell.Style = "Normal"
$ws.Range('E14').Value = '  +5.06%  '

$ws.Range('E15').Value = '  +0.32%  '

$dCell = $ws.Range('D16')
$dCell.Value = '''3.625.81'
$dCell.Style = "Normal"
$ws.Range('E16').Value = '  +3.91%  '

$dCell = $ws.Range('D17')
$dCell.Value = '''7.16'
$dCell.Style = "Normal"
$ws.Range('E17').Value = '  -0.31%  '

$dCell = $ws.Range('D18')
$dCell.Value = '''62.958.74'
$dCell.Style = "Normal"
$ws.Range('E18').Value = '  +5.91%  '

$dCell = $ws.Range('D19')
$dCell.Value = '''3.112.64'
$dCell.Style = "Normal"
$ws.Range('E19').Value = '  +4.14%  '

$dCell = $ws.Range('D20')
$dCell.Value = '''465.67'
$dCell.Style = "Normal"
$ws.Range('E20').Value = '  +7.09%  '

$dCell = $ws.Range('D21')
$dCell.Value = '''14.14'
$dCell.Style = "Normal"
$ws.Range('E21').Value = '  +4.37%  '

$ws.Range('E22').Value = '  +1.13%  '

$ws.Range('E23').Value = '  +7.32%  '

$dCell = $ws.Range('D24')
$dCell.Value = '''13.28'
$dCell.Style = "Normal"
$ws.Range('E24').Value = '  -1.15%  '

$dCell = $ws.Range('D25')
$dCell.Value = '''81.88'
$dCell.Style = "Normal"
$ws.Range('E25').Value = '  +2.28%  '

$dCell = $ws.Range('D26')
$dCell.Value = '''0.999'
$dCell.Style = "Normal"
$ws.Range('E26').Value = '  -0.02%  '

$dCell = $ws.Range('D27')
$dCell.Value = '''8.38'
$dCell.Style = "Normal"
$ws.Range('E27').Value = '  +7.79%  '

$ws.Range('E28').Value = '  +0.28%  '

$dCell = $ws.Range('D29')
$dCell.Value = '''2.67'
$dCell.Style = "Normal"
$ws.Range('E29').Value = '  +5.14%  '

$ws.Range('E30').Value = '  -0.13%  '

$dCell = $ws.Range('D31')
$dCell.Value = '''6.82'
$dCell.Style = "Normal"
$ws.Range('E31').Value = '  +8.97%  '

$dCell = $ws.Range('D32')
$dCell.Value = '''26.90'
$dCell.Style = "Normal"
$ws.Range('E32').Value = '  +4.48%  '

$dCell = $ws.Range('D34')
$dCell.Value = '''0.0₃0864'
$dCell.Style = "Normal"
$ws.Range('E34').Value = '  +11.03%  '

$ws.Range('E35').Value = '  +15.46%  '

$ws.Range('E36').Value = '  +4.35%  '

$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$dCell = $ws.Range('D37')
$dCell.Value = '''6.02'
$dCell.Style = "Normal"
$ws.Range('E37').Value = '  +2.21%  '

$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$dCell = $ws.Range('D38')
$dCell.Value = '''3.29'
$dCell.Style = "Normal"
$ws.Range('E38').Value = '  +18.81%  '

$dCell = $ws.Range('D39')
$dCell.Value = '''50.85'
$dCell.Style = "Normal"
$ws.Range('E39').Value = '  +4.11%  '

$dCell = $ws.Range('D40')
$dCell.Value = '''438.41'
$dCell.Style = "Normal"
$ws.Range('E40').Value = '  +9.44%  '

$dCell = $ws.Range('D41')
$dCell.Value = '''8.71'
$dCell.Style = "Normal"
$ws.Range('E41').Value = '  +1.25%  '

$dCell = $ws.Range('D42')
$dCell.Value = '''2.923.66'
$dCell.Style = "Normal"
$ws.Range('E42').Value = '  +5.94%  '

$ws.Range('E43').Value = '  +4.34%  '

$dCell = $ws.Range('D44')
$dCell.Value = '''0.279'
$dCell.Style = "Normal"
$ws.Range('E44').Value = '  +11.06%  '

$ws.Range('E45').Value = '  +5.54%  '

$dCell = $ws.Range('D46')
$dCell.Value = '''2.16'
$dCell.Style = "Normal"
$ws.Range('E46').Value = '  +7.62%  '

$dCell = $ws.Range('D47')
$dCell.Value = '''35.31'
$dCell.Style = "Normal"
$ws.Range('E47').Value = '  +2.18%  '

$dCell = $ws.Range('D49')
$dCell.Value = '''123.15'
$dCell.Style = "Normal"
$ws.Range('E49').Value = '  +0.13%  '

$ws.Range('E50').Value = '  +0.61%  '

$dCell = $ws.Range('D51')
$dCell.Value = '''24.48'
$dCell.Style = "Normal"
$ws.Range('E51').Value = '  +4.03%  '
